# The sheet originally has 13 columns (A:M). Katalon AI regenerated the
# report so that only the "input_Name" column (originally column I) is
# kept; every other column is removed and column I's data/formatting
# (including its header style and its already-blank row-2 cell) slides
# into column A.
#
# Delete from the right-hand side first (J:M) so the left-hand delete
# (A:H) doesn't disturb the as-yet-undeleted J:M addresses.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns("J:M").Delete()
$ws.Columns("A:H").Delete()
